# Version 1.1 of the pitch deck: add six new "Title and Content" slides
# after the existing "Solution" slide, fill in the Solution slide's body
# placeholder, and give each new slide its title text.

$p = $ppt.ActivePresentation

# The existing slides (1-5) all use the "Title and Content" custom layout
# (slide master layout #2). Re-use that same layout for the new slides so
# they match the rest of the deck.
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

# --- Slide 5 ("Solution") : fill in the, so far empty, body placeholder ---
$solutionSlide = $p.Slides.Item(5)
$solutionSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Solution justification"

# --- New slides appended at the end of the deck -----------------------
$newTitles = @(
    "Market limitation",
    "Target market ",
    "Business Model/Money making points",
    "Team members selections and contribution",
    "Starting Strategy"
)

foreach ($title in $newTitles) {
    $idx = $p.Slides.Count + 1
    $slide = $p.Slides.AddSlide($idx, $titleAndContent)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $title
}

# Last new slide has a title typed in two separate passes ("Tools used "
# then "for the team" appended afterwards), which PowerPoint records as
# two separate runs.
$lastIdx = $p.Slides.Count + 1
$toolsSlide = $p.Slides.AddSlide($lastIdx, $titleAndContent)
$toolsTitle = $toolsSlide.Shapes.Item(1).TextFrame.TextRange
$toolsTitle.Text = "Tools used "
$toolsTitle.InsertAfter("for the team") | Out-Null
